$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.796.53"
$ws.Range("E2").Value = "  -3.23%  "
$ws.Range("D3").Value = "1.613.43"
$ws.Range("E3").Value = "  -3.60%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.26%  "
$ws.Range("E5").Value = "  +0.16%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "306.73"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.00%  "
$ws.Range("E7").Value = "  -0.54%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3795"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.90%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.003"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.40%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.348"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.05%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "48.55"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -5.67%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08411"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.65%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "23.74"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -6.30%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.013"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.38%  "
$ws.Range("E15").Value = "  -4.13%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.440"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.90%  "
$ws.Range("D17").Value = "1.613.78"
$ws.Range("E17").Value = "  -3.31%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "92.94"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.55%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06912"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.41%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.94"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -5.35%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.780"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.00%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.001"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.06%  "
$ws.Range("E23").Value = "  -4.08%  "
$ws.Range("D24").Value = "23.797.31"
$ws.Range("E24").Value = "  -3.20%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.433"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.95%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.796"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.74%  "
$ws.Range("E27").Value = "  -4.85%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "156.91"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.20%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "139.03"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -5.16%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.251"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -10.48%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.723"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -7.81%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.483"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.93%  "
$ws.Range("D33").Value = "1.789.50"
$ws.Range("E33").Value = "  -3.55%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08106"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.99%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9633"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.64%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02860"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -8.25%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.513"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -6.79%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2639"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -6.07%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.09127"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.77%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "10.35"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.01%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "13.40"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.45%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.423"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -6.19%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.7407"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -6.43%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "15.76"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.10%  "
$ws.Range("E45").Value = "  -4.64%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.429"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.82%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.050"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.30%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.001"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.11%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.08209"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "132.23"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.75%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.190"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -10.18%  "
